$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra empty template rows (3-10), keep header + first data row
$ws.Range("3:10").Delete()

# Insert a new "status" column before the existing "remark" column (K)
$ws.Range("K1").EntireColumn.Insert()

# New column header
$ws.Range("K1").Value = "สถานะ"

# Update the sample data row with the new values
$ws.Range("A2").Value = "สกลราชวิทยานุกูล"
$ws.Range("I2").Value = 102800
$ws.Range("J2").Value = "สพฐ"

# Update / extend the comment on A1 with guidance about the new status column
[void]$ws.Range("A1").Comment.Text("Imported Author:`nOptional: ใส่หรือไม่ใส่ก็ได้`nValue: 1 สถานะใช้งาน, 0 ระงับการใช้งาน`nDefault: 1`n`t-Nuttasak Tawan")
